# Insert a new "editor_name" column before the existing "coverImage" column (S),
# which pushes coverImage (header + data) from column S to column T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S; existing column S ("coverImage") and everything in
# it slides over to column T.
$ws.Columns("S:S").Insert()

# New header cell, matching the bold/centered header style used by the rest
# of row 1 (copy formatting from the neighboring header cell).
$ws.Range("S1").Value = "editor_name"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

# Fill the new column's data rows (2-17) with the same empty-placeholder text
# ("''") already used by the neighboring editor_img/editor_bio columns on
# this sheet. Copy the literal values from those existing cells instead of
# re-typing them, so Excel's leading-quote "treat as text" autocorrect
# doesn't swallow a character.
$ws.Range("Q2:Q17").Copy()
$ws.Range("S2:S17").PasteSpecial(-4163)
